$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 14 timestamp (recalculated value from scheduled task)
$ws.Range("A14").Value = 45865.62522155092

# Copy the date/time number format from A14 so the new timestamp cell matches styling
$ws.Range("A15").NumberFormat = $ws.Range("A14").NumberFormat

# Append new row 15 with the latest sensor reading
$ws.Range("A15").Value = 45865.66692335347
$ws.Range("B15").Value = 2025
$ws.Range("C15").Value = 30
$ws.Range("D15").Value = 18.72
$ws.Range("E15").Value = 76.31
$ws.Range("F15").Value = 253.13
$ws.Range("G15").Value = 10.75
$ws.Range("H15").Value = "ESE"
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = "16:00:22"
